# Apply "suspension sample" CI-compliance edit:
#  - Rename the "section_thickness_unit list" sheet to "suspension_entity list"
#    and replace its values (um/mm/cm -> cell/nuclei).
#  - Rename the "area_unit list" sheet to "suspension_enriched list"
#    and replace its values (mm^2/um^2 -> yes/no).
#  - Replace the section/area related columns (M:R) on the "Export as TSV"
#    sheet with the new suspension related columns (M:Q), dropping the
#    now-unused last column.
#  - Update the column header comments and data validations to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the two lookup-list sheets first so that any data validation
#    formulas we create below can refer to them by their new names.
# ---------------------------------------------------------------------
$wsEntity = $wb.Worksheets.Item("section_thickness_unit list")
$wsEntity.Name = "suspension_entity list"

$wsEnriched = $wb.Worksheets.Item("area_unit list")
$wsEnriched.Name = "suspension_enriched list"

# ---------------------------------------------------------------------
# 2. Update the contents of the renamed lookup-list sheets.
# ---------------------------------------------------------------------
$wsEntity.Range("A1").Value = "cell"
$wsEntity.Range("A2").Value = "nuclei"
$wsEntity.Range("A3").Clear()

$wsEnriched.Range("A1").Value = "yes"
$wsEnriched.Range("A2").Value = "no"

# ---------------------------------------------------------------------
# 3. Update the main "Export as TSV" sheet.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Export as TSV")

# Preserve the existing "notes" comment text (currently on R1) so it can be
# reapplied to its new home at Q1.
$notesCommentText = $ws.Range("R1").Comment.Text()

# --- Column M: section_thickness_value -> suspension_entity ---
$ws.Range("M1").Value = "suspension_entity"
[void]$ws.Range("M1").Comment.Text("The type of single cell entity derived from isolation protocol.")

$mValidation = $ws.Range("M2:M1048576").Validation
$mValidation.Delete()
$mValidation.Add(3, 1, 1, "='suspension_entity list'!`$A`$1:`$A`$2")
$mValidation.ErrorTitle = "Value must come from list"
$mValidation.ErrorMessage = "Value must be one of: cell / nuclei."
$mValidation.IgnoreBlank = $true
$mValidation.InCellDropdown = $true
$mValidation.ShowInput = $true
$mValidation.ShowError = $true

# --- Column N: section_thickness_unit -> suspension_entity_number ---
$ws.Range("N1").Value = "suspension_entity_number"
[void]$ws.Range("N1").Comment.Text("Total number of cell/nuclei yielded post dissociation and enrichment.")

$nValidation = $ws.Range("N2:N1048576").Validation
$nValidation.Delete()
$nValidation.Add(1, 1, 1, "-2147483647", "2147483647")
$nValidation.ErrorTitle = "Not an integer"
$nValidation.ErrorMessage = "The values in this column must be integers."
$nValidation.IgnoreBlank = $true
$nValidation.ShowInput = $true
$nValidation.ShowError = $true

# --- Column O: section_index_number -> suspension_enriched ---
$ws.Range("O1").Value = "suspension_enriched"
[void]$ws.Range("O1").Comment.Text("Was the cell/nuclei population enriched?")

$oValidation = $ws.Range("O2:O1048576").Validation
$oValidation.Delete()
$oValidation.Add(3, 1, 1, "='suspension_enriched list'!`$A`$1:`$A`$2")
$oValidation.ErrorTitle = "Value must come from list"
$oValidation.ErrorMessage = "Value must be one of: yes / no."
$oValidation.IgnoreBlank = $true
$oValidation.InCellDropdown = $true
$oValidation.ShowInput = $true
$oValidation.ShowError = $true

# --- Column P: area_value -> suspension_enriched_target (free text, no validation) ---
$ws.Range("P1").Value = "suspension_enriched_target"
[void]$ws.Range("P1").Comment.Text("If the suspension was enriched, then this is the target of the enrichment.")
$ws.Range("P2:P1048576").Validation.Delete()

# --- Column Q: area_unit -> notes (free text, no validation) ---
$ws.Range("Q1").Value = "notes"
[void]$ws.Range("Q1").Comment.Text($notesCommentText)
$ws.Range("Q2:Q1048576").Validation.Delete()

# --- Column R: notes column is no longer needed; drop it entirely ---
$ws.Range("R1").Clear()
